# correção nos dados e inicio da analise PNAD 2009
#
# The header row (row 2) had two placeholder "unnamed" labels (columns B
# and F) that should read "total", matching column C. Columns C/D/E keep
# their existing text.
#
# The data block had two label-only "divider" rows ("situação do
# domicílio" above the urbana/rural split, and "grandes regiões e
# unidades da federação" above the region list) that carried no figures.
# Those rows are removed entirely, so every row below slides up and the
# last two (previously unused/overflow) rows at the bottom disappear.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix the second header row -------------------------------------------------
$ws.Range("B2").Value = "total"
$ws.Range("F2").Value = "total"

# --- drop the two label-only divider rows --------------------------------------
# "situação do domicílio" (row 5) sits right above urbana/rural.
$ws.Rows(5).Delete()

# After that deletion, "grandes regiões e unidades da federação" (was row
# 8) is now row 7, sitting right above "norte".
$ws.Rows(7).Delete()
